$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: human-readable headers (translated / accented Spanish labels) ---
$ws.Range("A1").Value = "Número de edificios"
$ws.Range("B1").Value = "Comarca nombre"
$ws.Range("C1").Value = "Comarca código"
$ws.Range("D1").Value = "Plantas sobre rasante"
$ws.Range("E1").Value = "Provincia código"
$ws.Range("F1").Value = "Aragón"
$ws.Range("G1").Value = "Municipio código"
$ws.Range("H1").Value = "Provincia nombre"
$ws.Range("I1").Value = "Municipio nombre"

# --- Row 2: SDMX/IAEST dimension & measure concept identifiers ---
$ws.Range("A2").Value = "iaest-measure:numero-de-edificios"
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-dimension:plantas-sobre-rasante"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = "sdmx-dimension:refArea"
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# --- Row 3: role (measure / dimension) ---
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "dim"
$ws.Range("I3").Value = "dim"

# --- Row 4: datatype / codelist concept reference ---
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "URI-comarca"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "skos:Concept"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = "URI-Provincia"
$ws.Range("I4").Value = "URI-Municipio"

# --- Row 5: mapping file reference moved from G5 to D5 ---
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats - carry over the existing cell style
$ws.Range("D5").Value = "mapping-plantas-sobre-rasante.xlsx"
$ws.Range("G5").Clear()
